# Updates the cryptos list with refreshed prices / 1h volume percentages,
# as captured by the scheduled GitHub Actions run on
# Sun May 12 19:34:33 UTC 2024.
#
# All Price (D) and Volume(1h) (E) cells in this sheet are stored as plain
# text, even when their content happens to look like a number (e.g.
# "597.51"). Assigning such strings straight to Range.Value makes Excel
# "helpfully" coerce them into real numbers, which would change the
# underlying cell type. To keep these cells as text - exactly like the
# source data - we stage the text in a scratch cell that has been
# explicitly formatted as Text ("@"), copy it, and paste just the value
# into the destination with PasteSpecial. That preserves the text cell
# type without leaving the destination cell's style touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Range("ZZ1")

function Set-TextValue($cellRef, $val) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $val
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "61.336.24"
$ws.Range("E2").Value = "  +0.40%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.927.24"
$ws.Range("E3").Value = "  -0.03%  "

# Row 5 - BNB
Set-TextValue "D5" "597.51"
$ws.Range("E5").Value = "  +0.57%  "

# Row 6 - Solana
Set-TextValue "D6" "145.11"
$ws.Range("E6").Value = "  -0.75%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.01%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -0.99%  "

# Row 9 - Toncoin
Set-TextValue "D9" "6.99"
$ws.Range("E9").Value = "  +1.42%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -2.62%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  -0.62%  "

# Row 12 - ShibaInu
$ws.Range("E12").Value = "  -1.13%  "

# Row 13 - Avalanche
Set-TextValue "D13" "33.45"
$ws.Range("E13").Value = "  -0.97%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +0.25%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "3.411.08"
$ws.Range("E15").Value = "  -0.01%  "

# Row 16 - WrappedBTC
Set-TextValue "D16" "61.292.34"
$ws.Range("E16").Value = "  +0.42%  "

# Row 17 & 18 - Polkadot and WrappedEther swapped rank positions
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D17" "6.69"
$ws.Range("E17").Value = "  -0.51%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D18" "2.921.44"
$ws.Range("E18").Value = "  -0.23%  "

# Row 19 - BitcoinCash
Set-TextValue "D19" "431.32"
$ws.Range("E19").Value = "  -0.23%  "

# Row 20 - Chainlink
Set-TextValue "D20" "13.47"
$ws.Range("E20").Value = "  -0.06%  "

# Row 21 - Polygon
Set-TextValue "D21" "0.675"
$ws.Range("E21").Value = "  -1.20%  "

# Row 22 - Uniswap
Set-TextValue "D22" "7.07"
$ws.Range("E22").Value = "  -0.37%  "

# Row 23 - Litecoin
Set-TextValue "D23" "81.77"
$ws.Range("E23").Value = "  +0.47%  "

# Row 24 - RenderToken
Set-TextValue "D24" "10.84"
$ws.Range("E24").Value = "  -2.19%  "

# Row 25 - Fetch.AI
$ws.Range("E25").Value = "  -2.03%  "

# Row 26 - InternetComputer(DFINITY)
Set-TextValue "D26" "11.71"
$ws.Range("E26").Value = "  -2.48%  "

# Row 27
$ws.Range("E27").Value = "  +0.06%  "

# Row 28
Set-TextValue "D28" "2.22"
$ws.Range("E28").Value = "  -4.56%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  -0.71%  "

# Row 30 - NEARProtocol
$ws.Range("E30").Value = "  -2.70%  "

# Row 31 & 32 - Hedera and EthereumClassic swapped rank positions
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D31" "26.62"
$ws.Range("E31").Value = "  +0.43%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D32" "0.109"
$ws.Range("E32").Value = "  +1.23%  "

# Row 33 - FirstDigitalUSD
$ws.Range("E33").Value = "  -0.03%  "

# Row 34 - PEPE
$ws.Range("D34").Value = "0.0₃0883"
$ws.Range("E34").Value = "  +3.27%  "

# Row 35 - Mantle
$ws.Range("E35").Value = "  -0.26%  "

# Row 36 - Filecoin
Set-TextValue "D36" "5.62"
$ws.Range("E36").Value = "  -0.40%  "

# Row 37 - dogwifhat
$ws.Range("E37").Value = "  -2.91%  "

# Row 38 - Stacks
$ws.Range("E38").Value = "  -0.25%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  -1.89%  "

# Row 40 - Cosmos
Set-TextValue "D40" "8.56"
$ws.Range("E40").Value = "  -0.49%  "

# Row 41 - Arweave
Set-TextValue "D41" "42.17"
$ws.Range("E41").Value = "  +5.96%  "

# Row 42 - TheGraph
Set-TextValue "D42" "0.281"
$ws.Range("E42").Value = "  -2.01%  "

# Row 43 - VeChain
Set-TextValue "D43" "0.0345"
$ws.Range("E43").Value = "  -0.37%  "

# Row 44 - Maker
Set-TextValue "D44" "2.699.23"
$ws.Range("E44").Value = "  -1.06%  "

# Row 45 - Monero
Set-TextValue "D45" "133.63"
$ws.Range("E45").Value = "  +2.25%  "

# Row 46 - Bittensor
Set-TextValue "D46" "361.19"
$ws.Range("E46").Value = "  -3.90%  "

# Row 47 - USDe
$ws.Range("E47").Value = "  +0.09%  "

# Row 48 - InjectiveProtocol
Set-TextValue "D48" "23.57"
$ws.Range("E48").Value = "  -2.47%  "

# Row 49 - Stellar
$ws.Range("E49").Value = "  -1.19%  "

# Row 50 - ThetaToken
$ws.Range("E50").Value = "  -1.75%  "

# Row 51 - Cronos
$ws.Range("E51").Value = "  -2.51%  "
